$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.986
$ws.Range("A21").Value = -20.269
$ws.Range("A23").Value = -20.596
$ws.Range("C24").Value = -12.218
$ws.Range("A25").Value = -21.731
$ws.Range("C28").Value = -13.179
$ws.Range("C36").Value = -12.732
$ws.Range("C45").Value = -13.335
$ws.Range("C48").Value = -11.205
$ws.Range("C49").Value = -13.367
$ws.Range("C52").Value = -11.655
$ws.Range("A53").Value = -21.836
$ws.Range("C53").Value = -12.789
$ws.Range("C54").Value = -13.34
$ws.Range("A57").Value = -22.121
$ws.Range("A59").Value = -22.37
$ws.Range("A69").Value = -21.519
$ws.Range("C70").Value = -11.492
$ws.Range("A79").Value = -21.15
$ws.Range("A83").Value = -21.967
$ws.Range("C86").Value = -13.9
$ws.Range("C87").Value = -13.313
$ws.Range("A93").Value = -21.508
$ws.Range("C101").Value = -12.721

$wb.Save()
